$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170378684997559
$ws.Range("B1").Value = 2.320566654205322
$ws.Range("C1").Value = 3.227429151535034
$ws.Range("D1").Value = 1.429557681083679
$ws.Range("E1").Value = 1.156444549560547
